$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new attendance day (2019/11/30) has been recorded: add a new column I
# with that date as header, mirroring the style of the existing "MAC"
# attendance columns (E:H), and mark every student present (value 1).

# -- Header cell I1: "2019/11/30" --
$ws.Range("I1").NumberFormat = "@"          # force text so the date-like string isn't auto-converted
$ws.Range("I1").Value = "2019/11/30"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)         # xlPasteFormats: reuse H1's style (centered)

# -- Data cells I2:I8: attendance value 1, styled like column H --
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 9).PasteSpecial(-4122)
}

# -- I9: keep same (unstyled) look as H9 --
$ws.Cells.Item(9, 9).Value = 1
